# Updated cryptos list (price/volume refresh + MXToken/TheSandbox row swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.620.78'
$ws.Range('E2').Value = '  +2.67%  '
$ws.Range('D3').Value = '1.855.18'
$ws.Range('E3').Value = '  +2.46%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.034'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +2.75%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '322.17'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.36%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.030'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.51%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4398'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.67%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3792'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07436'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8795'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.21%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.68'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.61%  '
$ws.Range('D12').Value = '1.859.46'
$ws.Range('E12').Value = '  -8.40%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.533'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.87%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.710'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.37%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.07215'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.63%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '83.29'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.34%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.036'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009076'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.83%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.030'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.53%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.47'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.12%  '
$ws.Range('D21').Value = '27.647.52'
$ws.Range('E21').Value = '  +2.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.281'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.99%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.41'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.99%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '158.19'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.92%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.926'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.21%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '18.79'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.96%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.982'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.97%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.291'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.68%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '117.52'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.42%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09069'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.67%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.206'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.34%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7647'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.16%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.540'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.92%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.891'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.11%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.031'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.80%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.153'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.36%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01981'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.27%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05324'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.23%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.831'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.83%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5173'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.95%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1681'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.39%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.797'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.539'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.75%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '109.17'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.28%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.54'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.17%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.717'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.23%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4661'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.89%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.06410'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.854'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.24%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '39.44'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.56%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '64.21'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.02%  '
